$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: literal date-like text strings (must NOT be auto-converted to date serials) ---
# Strategy: write each date as a text-formula result, recalc, then Copy + PasteSpecial(xlPasteValues)
# over the whole A755:A774 block. PasteSpecial values-only does not re-run the "looks like a date"
# literal-input heuristic that a direct .Value = "2024-09-02" assignment would trigger, and it does
# not leave any NumberFormat/quotePrefix style behind either.
$ws.Range("A755").Formula = '="2024-09-02"'
$ws.Range("A756").Formula = '="2024-09-03"'
$ws.Range("A757").Formula = '="2024-09-04"'
$ws.Range("A758").Formula = '="2024-09-05"'
$ws.Range("A759").Formula = '="2024-09-06"'
$ws.Range("A760").Formula = '="2024-09-09"'
$ws.Range("A761").Formula = '="2024-09-10"'
$ws.Range("A762").Formula = '="2024-09-11"'
$ws.Range("A763").Formula = '="2024-09-12"'
$ws.Range("A764").Formula = '="2024-09-13"'
$ws.Range("A765").Formula = '="2024-09-16"'
$ws.Range("A766").Formula = '="2024-09-17"'
$ws.Range("A767").Formula = '="2024-09-18"'
$ws.Range("A768").Formula = '="2024-09-19"'
$ws.Range("A769").Formula = '="2024-09-20"'
$ws.Range("A770").Formula = '="2024-09-23"'
$ws.Range("A771").Formula = '="2024-09-24"'
$ws.Range("A772").Formula = '="2024-09-25"'
$ws.Range("A773").Formula = '="2024-09-26"'
$ws.Range("A774").Formula = '="2024-09-27"'

$excel.Calculate()

$dateRange = $ws.Range("A755:A774")
$dateRange.Copy()
$dateRange.PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = $false

# --- Numeric columns C..J: plain numeric literals, written directly ---
$ws.Range("C755").Value = 1885.400024414062
$ws.Range("D755").Value = 683.5999755859375
$ws.Range("E755").Value = 73.80999755859375
$ws.Range("F755").Value = 296.8999938964844
$ws.Range("G755").Value = 1303.849975585938
$ws.Range("H755").Value = 28487.57955932617
$ws.Range("I755").Value = 0
$ws.Range("J755").Value = 512.4776708814596

$ws.Range("C756").Value = 1901.949951171875
$ws.Range("D756").Value = 689.4000244140625
$ws.Range("E756").Value = 74.47000122070312
$ws.Range("F756").Value = 297.1499938964844
$ws.Range("G756").Value = 1320.25
$ws.Range("H756").Value = 28715.86001586914
$ws.Range("I756").Value = 0.008013332830455756
$ws.Range("J756").Value = 516.5843250264095

$ws.Range("C757").Value = 1901.300048828125
$ws.Range("D757").Value = 688.9500122070312
$ws.Range("E757").Value = 74.16000366210938
$ws.Range("F757").Value = 298.9500122070312
$ws.Range("G757").Value = 1327.099975585938
$ws.Range("H757").Value = 28721.03076171875
$ws.Range("I757").Value = 0.0001800658537390795
$ws.Range("J757").Value = 516.6773442239236

$ws.Range("C758").Value = 1879.449951171875
$ws.Range("D758").Value = 687.5
$ws.Range("E758").Value = 76
$ws.Range("F758").Value = 290.6000061035156
$ws.Range("G758").Value = 1312.349975585938
$ws.Range("H758").Value = 28602.74978637695
$ws.Range("I758").Value = -0.004118270556621158
$ws.Range("J758").Value = 514.5495271299329

$ws.Range("C759").Value = 1872.349975585938
$ws.Range("D759").Value = 673.5499877929688
$ws.Range("E759").Value = 74.72000122070312
$ws.Range("F759").Value = 283.6000061035156
$ws.Range("G759").Value = 1289.699951171875
$ws.Range("H759").Value = 28191.60983276367
$ws.Range("I759").Value = -0.01437414083205038
$ws.Range("J759").Value = 507.1533197619024

$ws.Range("C760").Value = 1892.400024414062
$ws.Range("D760").Value = 664.1500244140625
$ws.Range("E760").Value = 74.33999633789062
$ws.Range("F760").Value = 281.5499877929688
$ws.Range("G760").Value = 1237.150024414062
$ws.Range("H760").Value = 28036.46997070312
$ws.Range("I760").Value = -0.005503050836077006
$ws.Range("J760").Value = 504.3624292615674

$ws.Range("C761").Value = 1922.449951171875
$ws.Range("D761").Value = 664.5999755859375
$ws.Range("E761").Value = 78.05000305175781
$ws.Range("F761").Value = 285.75
$ws.Range("G761").Value = 1250.300048828125
$ws.Range("H761").Value = 28561.49984741211
$ws.Range("I761").Value = 0.01872667555001102
$ws.Range("J761").Value = 513.8074608339642

$ws.Range("C762").Value = 1957.599975585938
$ws.Range("D762").Value = 689.75
$ws.Range("E762").Value = 81.94999694824219
$ws.Range("F762").Value = 288.0499877929688
$ws.Range("G762").Value = 1237.699951171875
$ws.Range("H762").Value = 29297.64938354492
$ws.Range("I762").Value = 0.0257741904334731
$ws.Range("J762").Value = 527.0504321756381

$ws.Range("C763").Value = 1996.400024414062
$ws.Range("D763").Value = 729.1500244140625
$ws.Range("E763").Value = 81.69999694824219
$ws.Range("F763").Value = 291.7000122070312
$ws.Range("G763").Value = 1237.300048828125
$ws.Range("H763").Value = 29912.80038452148
$ws.Range("I763").Value = 0.02099659917843318
$ws.Range("J763").Value = 538.1166988468499

$ws.Range("C764").Value = 1988.050048828125
$ws.Range("D764").Value = 713.7000122070312
$ws.Range("E764").Value = 83.11000061035156
$ws.Range("F764").Value = 289.9500122070312
$ws.Range("G764").Value = 1241.5
$ws.Range("H764").Value = 29812.18057250977
$ws.Range("I764").Value = -0.003363771051799782
$ws.Range("J764").Value = 536.3065974727788

$ws.Range("C765").Value = 1989.900024414062
$ws.Range("D765").Value = 714.2000122070312
$ws.Range("E765").Value = 84.69999694824219
$ws.Range("F765").Value = 290.3999938964844
$ws.Range("G765").Value = 1226.599975585938
$ws.Range("H765").Value = 29926.49987792969
$ws.Range("I765").Value = 0.003834650911961043
$ws.Range("J765").Value = 538.3631460558686

$ws.Range("C766").Value = 2006.550048828125
$ws.Range("D766").Value = 731.0999755859375
$ws.Range("E766").Value = 82
$ws.Range("F766").Value = 284.2999877929688
$ws.Range("G766").Value = 1193.800048828125
$ws.Range("H766").Value = 29823.24993896484
$ws.Range("I766").Value = -0.003450117433913109
$ws.Range("J766").Value = 536.5057299798849

$ws.Range("C767").Value = 1987.800048828125
$ws.Range("D767").Value = 743.25
$ws.Range("E767").Value = 80.81999969482422
$ws.Range("F767").Value = 282.8500061035156
$ws.Range("G767").Value = 1166.400024414062
$ws.Range("H767").Value = 29685.31034851074
$ws.Range("I767").Value = -0.004625236710834788
$ws.Range("J767").Value = 534.0242639820087

$ws.Range("C768").Value = 1998.599975585938
$ws.Range("D768").Value = 735.9500122070312
$ws.Range("E768").Value = 80.97000122070312
$ws.Range("F768").Value = 272.7000122070312
$ws.Range("G768").Value = 1121.300048828125
$ws.Range("H768").Value = 29455.8603515625
$ws.Range("I768").Value = -0.0077294120982553
$ws.Range("J768").Value = 529.8965703752243

$ws.Range("C769").Value = 2048.10009765625
$ws.Range("D769").Value = 746.5
$ws.Range("E769").Value = 83.44999694824219
$ws.Range("F769").Value = 277.3500061035156
$ws.Range("G769").Value = 1149.400024414062
$ws.Range("H769").Value = 30118.95037841797
$ws.Range("I769").Value = 0.02251131078642199
$ws.Range("J769").Value = 541.8252367556001

$ws.Range("C770").Value = 2082.39990234375
$ws.Range("D770").Value = 773.9500122070312
$ws.Range("E770").Value = 82.88999938964844
$ws.Range("F770").Value = 286.2999877929688
$ws.Range("G770").Value = 1162.75
$ws.Range("H770").Value = 30664.31942749023
$ws.Range("I770").Value = 0.01810717313253569
$ws.Range("J770").Value = 551.636160125111

$ws.Range("C771").Value = 2068.14990234375
$ws.Range("D771").Value = 781.8499755859375
$ws.Range("E771").Value = 83.79000091552734
$ws.Range("F771").Value = 291.7999877929688
$ws.Range("G771").Value = 1141.199951171875
$ws.Range("H771").Value = 30770.6690826416
$ws.Range("I771").Value = 0.003468188994144963
$ws.Range("J771").Value = 553.5493385844294

$ws.Range("C772").Value = 2061.60009765625
$ws.Range("D772").Value = 775.8499755859375
$ws.Range("E772").Value = 82.95999908447266
$ws.Range("F772").Value = 289.8500061035156
$ws.Range("G772").Value = 1118.449951171875
$ws.Range("H772").Value = 30542.33015441895
$ws.Range("I772").Value = -0.007420668286718119
$ws.Range("J772").Value = 549.4416325624621

$ws.Range("C773").Value = 2022.050048828125
$ws.Range("D773").Value = 772.0999755859375
$ws.Range("E773").Value = 81.83999633789062
$ws.Range("F773").Value = 290.5
$ws.Range("G773").Value = 1123.650024414062
$ws.Range("H773").Value = 30238.56976318359
$ws.Range("I773").Value = -0.009945553914831306
$ws.Range("J773").Value = 543.9771311827592

$ws.Range("C774").Value = 2031.300048828125
$ws.Range("D774").Value = 761.75
$ws.Range("E774").Value = 81.08999633789062
$ws.Range("F774").Value = 293.4500122070312
$ws.Range("G774").Value = 1130.25
$ws.Range("H774").Value = 30174.3701171875
$ws.Range("I774").Value = -0.002123104581297322
$ws.Range("J774").Value = 542.8222108434242

